$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the 📘 emoji (used in A2, A3, A4, A6, A7) with ⚠️
$ws.Range("A2").Value = "⚠️"
$ws.Range("A3").Value = "⚠️"
$ws.Range("A4").Value = "⚠️"
$ws.Range("A6").Value = "⚠️"
$ws.Range("A7").Value = "⚠️"

# Replace the 📙 emoji (used in A5) with +3 (force text so Excel
# doesn't interpret "+3" as the number 3)
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "+3"
